# Applies the Betfair Back/Lay odds updates described in the commit diff.
# Each line sets a single cell to its new numeric value; all other cells
# (labels, unrelated odds) are left untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 1.67
$ws.Range("O2").Value = 1.29
$ws.Range("X2").Value = 15.5
$ws.Range("AC2").Value = 9

# Row 3
$ws.Range("F3").Value = 1.54
$ws.Range("I3").Value = 2.84
$ws.Range("J3").Value = 1.55
$ws.Range("V3").Value = 1.54

# Row 5
$ws.Range("F5").Value = 5.7
$ws.Range("G5").Value = 6.8
$ws.Range("H5").Value = 1.56
$ws.Range("I5").Value = 1.66
$ws.Range("V5").Value = 2.52
$ws.Range("AF5").Value = 980

# Row 6
$ws.Range("F6").Value = 1.34
$ws.Range("G6").Value = 1.42
$ws.Range("H6").Value = 9.6
$ws.Range("I6").Value = 14
$ws.Range("J6").Value = 5.1
$ws.Range("N6").Value = 3.85
$ws.Range("O6").Value = 1.18
$ws.Range("Q6").Value = 1.66
$ws.Range("S6").Value = 2.42
$ws.Range("W6").Value = 3.35

# Row 7
$ws.Range("F7").Value = 1.73
$ws.Range("G7").Value = 1.86
$ws.Range("I7").Value = 6
$ws.Range("N7").Value = 2.82
$ws.Range("P7").Value = 1.86
$ws.Range("R7").Value = 1.27
$ws.Range("V7").Value = 1.2
$ws.Range("W7").Value = 2.16

# Row 8
$ws.Range("AC8").Value = 11.5

# Row 9
$ws.Range("N9").Value = 3.2

# Row 11
$ws.Range("I11").Value = 1.64
$ws.Range("O11").Value = 1.39
$ws.Range("R11").Value = 1.3
$ws.Range("V11").Value = 2.56

# Row 12
$ws.Range("K12").Value = 3.3
$ws.Range("Y12").Value = 9.4
$ws.Range("AB12").Value = 9.6
$ws.Range("AN12").Value = 38

# Row 13
$ws.Range("G13").Value = 4.8
$ws.Range("T13").Value = 2
$ws.Range("X13").Value = 13
$ws.Range("AE13").Value = 25
$ws.Range("AG13").Value = 19

# Row 14
$ws.Range("N14").Value = 2.66
$ws.Range("O14").Value = 1.52
$ws.Range("V14").Value = 1.65

# Row 15
$ws.Range("L15").Value = 1.36
$ws.Range("M15").Value = 1.05
$ws.Range("N15").Value = 4
$ws.Range("R15").Value = 1.41
$ws.Range("S15").Value = 2.96
$ws.Range("T15").Value = 1.69
$ws.Range("U15").Value = 2.18
$ws.Range("V15").Value = 1.29
$ws.Range("X15").Value = 21
$ws.Range("Y15").Value = 19.5
$ws.Range("Z15").Value = 36
$ws.Range("AA15").Value = 95
$ws.Range("AB15").Value = 13
$ws.Range("AC15").Value = 10.5
$ws.Range("AD15").Value = 20
$ws.Range("AE15").Value = 55
$ws.Range("AF15").Value = 16.5
$ws.Range("AG15").Value = 13
$ws.Range("AH15").Value = 21
$ws.Range("AI15").Value = 60
$ws.Range("AJ15").Value = 30
$ws.Range("AK15").Value = 25
$ws.Range("AL15").Value = 40
$ws.Range("AN15").Value = 16
$ws.Range("AO15").Value = 55

# Row 16
$ws.Range("F16").Value = 2.32
$ws.Range("G16").Value = 2.34
$ws.Range("H16").Value = 3.55
$ws.Range("I16").Value = 3.6
$ws.Range("L16").Value = 1.4
$ws.Range("U16").Value = 2.16
$ws.Range("W16").Value = 1.74
$ws.Range("Z16").Value = 23
$ws.Range("AB16").Value = 9.800000000000001
$ws.Range("AC16").Value = 7.4

# Row 17
$ws.Range("N17").Value = 3.6

# Row 18
$ws.Range("F18").Value = 3.7
$ws.Range("T18").Value = 1.86
$ws.Range("U18").Value = 2

# Row 19
$ws.Range("N19").Value = 4.4
$ws.Range("Q19").Value = 1.82
$ws.Range("V19").Value = 1.89

# Row 20
$ws.Range("F20").Value = 1.15
$ws.Range("R20").Value = 1.85
$ws.Range("T20").Value = 2.3

# Row 21
$ws.Range("G21").Value = 5.5
$ws.Range("N21").Value = 1.34
$ws.Range("P21").Value = 1.34
$ws.Range("W21").Value = 1.22
